$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imageName = "0a5dbcdc-2e44-4579-a576-c93d5ee55485.png"

# Existing data occupies rows 1-720 (A1:B720). Append 25 new rows (721-745)
# that mirror the "remote api config" rows already present at the tail of
# the sheet (value 0 in column A, the shared image filename in column B),
# with the final row carrying a value of 126 in column A.
$startRow = 721
$endRow = 745

for ($r = $startRow; $r -le $endRow; $r++) {
    if ($r -eq $endRow) {
        $ws.Cells.Item($r, 1).Value = 126
    } else {
        $ws.Cells.Item($r, 1).Value = 0
    }
    $ws.Cells.Item($r, 2).Value = $imageName
}
